# Repull data / push all data: update the dSF (column F) values on Sheet1
# for the rows whose upstream source data changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F2"  = 1
    "F3"  = 1
    "F5"  = -1
    "F8"  = -5
    "F9"  = -4
    "F16" = -6
    "F17" = -1
    "F23" = 7
    "F26" = 6
    "F27" = -9
    "F29" = -11
    "F30" = 6
    "F33" = 2
    "F34" = -1
    "F44" = 4
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
